# Update LR-pairs TPM values per new analysis run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MuSCs -> Ihh -> Hhip -> ECs): receptor stats + derived edge stats recomputed.
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.046397
$ws.Range("N2").Value = 0.139191
$ws.Range("O2").Value = 0.01970724914037141
$ws.Range("P2").Value = 0.01970724914037141
$ws.Range("Q2").Value = 0.003155769283333333
$ws.Range("R2").Value = 0.02840192355
$ws.Range("S2").Value = 0.01970724914037141
$ws.Range("T2").Value = 0.01970724914037141

# Row 3 (MuSCs -> Ihh -> Hhip -> FAPs): only derived specificity columns change.
$ws.Range("O3").Value = 0.6598912010221247
$ws.Range("P3").Value = 0.6598912010221247
$ws.Range("S3").Value = 0.6598912010221247
$ws.Range("T3").Value = 0.6598912010221247

# Row 4 (MuSCs -> Ihh -> Hhip -> MuSCs): receptor stats + derived edge stats recomputed.
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.754325
$ws.Range("N4").Value = 2.262975
$ws.Range("O4").Value = 0.3204015498375038
$ws.Range("P4").Value = 0.3204015498375038
$ws.Range("Q4").Value = 0.05130667208333334
$ws.Range("R4").Value = 0.46176004875
$ws.Range("S4").Value = 0.3204015498375038
$ws.Range("T4").Value = 0.3204015498375038
